$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.548.63'
$ws.Range("E2").Value = '  +6.17%  '
$ws.Range("D3").Value = '2.518.87'
$ws.Range("E3").Value = '  +4.08%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'325.60"
$ws.Range("E5").Value = '  +2.60%  '
$ws.Range("D6").Value = "'105.50"
$ws.Range("E6").Value = '  +3.25%  '
$ws.Range("E7").Value = '  +2.15%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").Value = "'0.543"
$ws.Range("E9").Value = '  +3.37%  '
$ws.Range("D10").Value = "'37.39"
$ws.Range("E10").Value = '  +5.85%  '
$ws.Range("D11").Value = "'0.0821"
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("E12").Value = '  +1.00%  '
$ws.Range("D13").Value = "'18.51"
$ws.Range("E13").Value = '  +0.34%  '
$ws.Range("D14").Value = "'7.23"
$ws.Range("E14").Value = '  +4.94%  '
$ws.Range("D15").Value = '2.910.97'
$ws.Range("E15").Value = '  +4.08%  '
$ws.Range("D16").Value = '2.515.75'
$ws.Range("E16").Value = '  +3.90%  '
$ws.Range("D17").Value = "'0.851"
$ws.Range("E17").Value = '  +3.00%  '
$ws.Range("D18").Value = '47.374.81'
$ws.Range("E18").Value = '  +6.18%  '
$ws.Range("E19").Value = '  +4.55%  '
$ws.Range("D20").Value = "'6.60"
$ws.Range("E20").Value = '  +4.35%  '
$ws.Range("D21").Value = '0.0₃0944'
$ws.Range("E21").Value = '  +2.99%  '
$ws.Range("D22").Value = "'71.13"
$ws.Range("E22").Value = '  +3.82%  '
$ws.Range("D23").Value = "'253.09"
$ws.Range("E23").Value = '  +4.27%  '
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = '  +5.13%  '
$ws.Range("E25").Value = '  +3.58%  '
$ws.Range("D26").Value = "'26.59"
$ws.Range("E26").Value = '  +5.70%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").Value = "'10.06"
$ws.Range("E28").Value = '  +5.83%  '
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = '  -4.09%  '
$ws.Range("D30").Value = "'35.45"
$ws.Range("E30").Value = '  +6.40%  '
$ws.Range("E31").Value = '  +6.95%  '
$ws.Range("D32").Value = "'49.70"
$ws.Range("E32").Value = '  +2.49%  '
$ws.Range("D33").Value = "'19.90"
$ws.Range("E33").Value = '  +2.54%  '
$ws.Range("E34").Value = '  +3.00%  '
$ws.Range("D35").Value = "'0.0786"
$ws.Range("E35").Value = '  +3.27%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  +4.11%  '
$ws.Range("D38").Value = "'4.63"
$ws.Range("E38").Value = '  +4.35%  '
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = '  +4.87%  '
$ws.Range("D40").Value = "'123.28"
$ws.Range("E40").Value = '  -2.70%  '
$ws.Range("E41").Value = '  +2.41%  '
$ws.Range("E42").Value = '  +2.66%  '
$ws.Range("D43").Value = "'21.74"
$ws.Range("E43").Value = '  +3.95%  '
$ws.Range("E44").Value = '  +3.65%  '
$ws.Range("D45").Value = '1.984.92'
$ws.Range("E45").Value = '  +2.48%  '
$ws.Range("E46").Value = '  +4.14%  '
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("E48").Value = '  +2.80%  '
$ws.Range("D49").Value = "'9.17"
$ws.Range("E49").Value = '  +0.48%  '
$ws.Range("D50").Value = "'5.47"
$ws.Range("E50").Value = '  +18.93%  '
$ws.Range("D51").Value = "'79.95"
$ws.Range("E51").Value = '  +5.63%  '
